$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B20").Value = "Fosso"
$ws.Range("C20").Value = "YES"
$ws.Range("C21").Value = "YES"

$ws.Range("E20").Select()
$excel.ActiveWindow.ScrollRow = 4
